$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Price" column (D) holds numbers formatted/typed as plain text in the
# source data (e.g. "26.134.36", "1.00", "0.0790"). Force these cells to the
# Text number format before assignment so Excel does not silently reinterpret
# them as numbers (which would drop significant trailing zeros / thousands dots).
$priceCells = @('D2', 'D3', 'D5', 'D7', 'D9', 'D10', 'D11', 'D12', 'D13', 'D14', 'D15', 'D16', 'D18', 'D20', 'D21', 'D22', 'D23', 'D24', 'D26', 'D27', 'D28', 'D29', 'D32', 'D33', 'D36', 'D37', 'D38', 'D39', 'D42', 'D43', 'D45', 'D47', 'D50', 'D51')
foreach ($addr in $priceCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = '26.134.36'
$ws.Range("E2").Value = '  +1.42%  '
$ws.Range("D3").Value = '1.641.68'
$ws.Range("E3").Value = '  +0.44%  '
$ws.Range("E4").Value = '  -0.17%  '
$ws.Range("D5").Value = '216.49'
$ws.Range("E5").Value = '  +0.34%  '
$ws.Range("E6").Value = '  +0.72%  '
$ws.Range("D7").Value = '1.00'
$ws.Range("E8").Value = '  +0.67%  '
$ws.Range("D9").Value = '0.0635'
$ws.Range("E9").Value = '  +0.28%  '
$ws.Range("D10").Value = '19.72'
$ws.Range("E10").Value = '  +0.70%  '
$ws.Range("D11").Value = '0.0790'
$ws.Range("E11").Value = '  -0.15%  '
$ws.Range("D12").Value = '1.869.45'
$ws.Range("E12").Value = '  +0.48%  '
$ws.Range("D13").Value = '4.28'
$ws.Range("E13").Value = '  +0.72%  '
$ws.Range("D14").Value = '1.638.78'
$ws.Range("E14").Value = '  +0.19%  '
$ws.Range("D15").Value = '0.544'
$ws.Range("E15").Value = '  -2.87%  '
$ws.Range("D16").Value = '0.0₃0762'
$ws.Range("E16").Value = '  -0.05%  '
$ws.Range("E17").Value = '  +0.07%  '
$ws.Range("D18").Value = '26.132.60'
$ws.Range("E18").Value = '  +1.32%  '
$ws.Range("E19").Value = '  -0.19%  '
$ws.Range("B20").Value = 'BitcoinCash'
$ws.Range("C20").Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range("D20").Value = '194.87'
$ws.Range("E20").Value = '  +1.20%  '
$ws.Range("B21").Value = 'Uniswap'
$ws.Range("C21").Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range("D21").Value = '4.42'
$ws.Range("E21").Value = '  -0.86%  '
$ws.Range("D22").Value = '10.03'
$ws.Range("E22").Value = '  +0.55%  '
$ws.Range("D23").Value = '6.26'
$ws.Range("E23").Value = '  -0.46%  '
$ws.Range("D24").Value = '1.79'
$ws.Range("E24").Value = '  -2.18%  '
$ws.Range("E25").Value = '  -0.11%  '
$ws.Range("D26").Value = '142.51'
$ws.Range("E26").Value = '  +0.24%  '
$ws.Range("D27").Value = '0.124'
$ws.Range("E27").Value = '  +1.04%  '
$ws.Range("D28").Value = '6.92'
$ws.Range("E28").Value = '  +0.30%  '
$ws.Range("D29").Value = '15.58'
$ws.Range("E29").Value = '  +0.60%  '
$ws.Range("E30").Value = '  +0.64%  '
$ws.Range("E31").Value = '  +1.79%  '
$ws.Range("D32").Value = '3.35'
$ws.Range("E32").Value = '  +0.66%  '
$ws.Range("D33").Value = '3.23'
$ws.Range("E33").Value = '  +0.24%  '
$ws.Range("E34").Value = '  +1.43%  '
$ws.Range("E35").Value = '  +1.34%  '
$ws.Range("D36").Value = '0.910'
$ws.Range("E36").Value = '  +0.69%  '
$ws.Range("D37").Value = '1.133.86'
$ws.Range("E37").Value = '  +0.21%  '
$ws.Range("D38").Value = '0.552'
$ws.Range("E38").Value = '  +1.49%  '
$ws.Range("D39").Value = '2.50'
$ws.Range("E39").Value = '  -0.27%  '
$ws.Range("E40").Value = '  +1.23%  '
$ws.Range("D42").Value = '100.34'
$ws.Range("E42").Value = '  -0.39%  '
$ws.Range("D43").Value = '5.48'
$ws.Range("E43").Value = '  -1.28%  '
$ws.Range("E44").Value = '  -0.39%  '
$ws.Range("D45").Value = '1.778.60'
$ws.Range("E45").Value = '  +0.51%  '
$ws.Range("E46").Value = '  -0.74%  '
$ws.Range("D47").Value = '56.73'
$ws.Range("E47").Value = '  +2.31%  '
$ws.Range("E48").Value = '  +4.29%  '
$ws.Range("E49").Value = '  +2.29%  '
$ws.Range("B50").Value = 'EnergySwap'
$ws.Range("C50").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D50").Value = '7.70'
$ws.Range("E50").Value = '  +3.55%  '
$ws.Range("B51").Value = 'Mantle'
$ws.Range("C51").Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range("D51").Value = '0.417'
$ws.Range("E51").Value = '  +0.06%  '
